$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: propagate existing formatting to the new cells BEFORE their
# source cells are overwritten, so the shared style entries get reused
# instead of new ones being minted. ---

# G2 / G6 become highlighted (orange) like C8/C9 currently are.
$ws.Range("C8").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C8").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# D8 becomes highlighted (yellow) like C7 currently is (SPEAKER moves there).
$ws.Range("C7").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# C8 / C9 become plain bordered (no fill, no center) like E8/E9 already are.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E9").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Step 2: move / update the cell values. ---

# SPEAKER moves from C7 to D8.
$ws.Range("D8").Value = "SPEAKER"

# TIMER0A (was in C8) now only lives in G2; TIMER2A is new in G6.
$ws.Range("G2").Value = "TIMER0A"
$ws.Range("G6").Value = "TIMER2A"

# C8 (TIMER0A) and C9 (TIMER0B) are cleared out.
$ws.Range("C8").Value = ""
$ws.Range("C9").Value = ""

# --- Step 3: C7 gets a new value (TRIMPOT) with a fresh highlight style:
# red fill + thin left/right border only (no top/bottom). ---
$ws.Range("C7").Borders.LineStyle = -4142          # clear existing borders
$ws.Range("C7").Borders.Item(7).LineStyle = 1      # xlEdgeLeft, thin
$ws.Range("C7").Borders.Item(10).LineStyle = 1     # xlEdgeRight, thin
$ws.Range("C7").Interior.Color = 255               # RGB(255,0,0) red
$ws.Range("C7").Value = "TRIMPOT"

# --- Step 4: selection moves to M7. ---
$ws.Range("M7").Select()
